$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1608416003758718
$ws.Range("C2").Value = 0.8594908078538727
$ws.Range("D2").Value = 1.64061426565291
$ws.Range("E2").Value = 1.280864655478053
$ws.Range("F2").Value = 1.318694559414439
$ws.Range("B3").Value = 0.05080766432823666
$ws.Range("C3").Value = 0.8009667910890287
$ws.Range("D3").Value = 1.083142076586474
$ws.Range("E3").Value = 1.040741118908288
$ws.Range("F3").Value = 1.081946107399986
$ws.Range("B4").Value = 0.05354799877553759
$ws.Range("C4").Value = 1.000578692683312
$ws.Range("D4").Value = 1.618272605442017
$ws.Range("E4").Value = 1.272113440476917
$ws.Range("F4").Value = 1.327501501701935
$ws.Range("B5").Value = 0.0459864616832071
$ws.Range("C5").Value = 0.7358968996768105
$ws.Range("D5").Value = 0.7047535131209822
$ws.Range("E5").Value = 0.8394959875550223
$ws.Range("F5").Value = 0.8791488123799778
$ws.Range("B6").Value = 0.1791212689832662
$ws.Range("C6").Value = 1.046431744547104
$ws.Range("D6").Value = 1.85744929065894
$ws.Range("E6").Value = 1.36288271346398
$ws.Range("F6").Value = 1.424142963195278
$ws.Range("B7").Value = 0.2249165586779437
$ws.Range("C7").Value = 0.9885917312309039
$ws.Range("D7").Value = 1.868591037025346
$ws.Range("E7").Value = 1.366964168157069
$ws.Range("F7").Value = 1.430123779954045
$ws.Range("B8").Value = 0.5374424523978547
$ws.Range("C8").Value = 1.313921765979487
$ws.Range("D8").Value = 2.669240729549001
$ws.Range("E8").Value = 1.633781114332333
$ws.Range("F8").Value = 1.690111122941772
$ws.Range("B9").Value = -0.2811530848179681
$ws.Range("C9").Value = 0.5494525610246418
$ws.Range("D9").Value = 0.4703157481287797
$ws.Range("E9").Value = 0.685795704367401
$ws.Range("F9").Value = 0.7660959708412387
$ws.Range("B10").Value = 0.1212484254788393
$ws.Range("C10").Value = 0.1212484254788393
$ws.Range("D10").Value = 0.01470118068109764
$ws.Range("E10").Value = 0.1212484254788393
